$d = $word.ActiveDocument

function Find-Text($searchText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "NOT FOUND:" $searchText
    }
    return $rng
}

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "REPLACE FAILED:" $old
    }
}

function Protect-Boundary($rng) {
    # Force a run split boundary by toggling a formatting property on/off across the range.
    $rng.Bold = 1
    $rng.Bold = 0
}

# Given the text that starts a run of concatenated "pieces" (already present verbatim,
# contiguous, in the document, e.g. after several Replace-Text calls merged them into one
# run), re-split that text back into separate <w:r> elements at each piece boundary.
function Split-Pieces($anchorText, $pieces) {
    $rng = Find-Text $anchorText
    $pos = $rng.Start
    for ($i = 0; $i -lt $pieces.Length - 1; $i++) {
        $pos = $pos + $pieces[$i].Length
        $boundaryEnd = $pos + $pieces[$i + 1].Length
        $bRng = $d.Range($pos, $boundaryEnd)
        Protect-Boundary $bRng
    }
}

# ---------------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------------
Replace-Text "Unraveling the Genetic Enigma of Disease" "Exploring the Human Body: An Introduction to Biology"

# ---------------------------------------------------------------------------
# Author name: "Helen Walsh" -> "Dr" + "." + " Emily Carter" (3 runs)
# ---------------------------------------------------------------------------
$authorRng = Find-Text "Helen Walsh"
$authorRng.Text = "Dr. Emily Carter"
Split-Pieces "Dr. Emily Carter" @("Dr", ".", " Emily Carter")

Write-Host "STAGE-AUTHOR OK"

# ---------------------------------------------------------------------------
# Email: "helen" + "." + "walsh@virology" (3 runs) -> "emilycarter@biostudies" (1 run)
# The trailing "." and "edu" runs are left as-is (but still need protecting from
# the automatic run-merge that happens after any edit in the same paragraph).
# ---------------------------------------------------------------------------
$helenRng = Find-Text "helen"
$emailStart = $helenRng.Start

$delRng = $d.Range($emailStart + 5, $emailStart + 20)
$delRng.Delete()

$helenRng2 = $d.Range($emailStart, $emailStart + 5)
$helenRng2.Text = "emilycarter@biostudies"

Split-Pieces "emilycarter@biostudies" @("emilycarter@biostudies", ".", "edu")

Write-Host "STAGE-EMAIL OK"

# ---------------------------------------------------------------------------
# Main body paragraph, segment A (sentences 1-3, separated by "." runs)
# ---------------------------------------------------------------------------
$a1old = "Within every cell, sequences of nucleotides arrange themselves, revealing the blueprint of life: genes, the microscopic directives that govern the intricate symphony of our biology"
$a1new = "Within the intricate tapestry of life, Biology unveils an astonishing symphony of interconnected systems that comprise the human body, an enigmatic marvel whose intricate workings continue to captivate and inspire awe"
$a2old = " The genetic code, a beacon of identity encoded within DNA's double helix, influences our development and guides our response to environmental cues"
$a2new = " Every cell, tissue, and organ collaborates in a delicate dance, performing a mesmerizing array of functions that sustain our existence"
$a3old = " While mysterious maladies may plague humanity, the decoding of our genetic inheritance can serve as a beacon of hope, illuminating pathways to healing"
$a3new = " This essay embarks on a journey into the realm of Biology, unraveling the enigmatic tapestry of the human body, and elucidating the profound mysteries that lie beneath the surface"

Replace-Text $a1old $a1new
Replace-Text $a2old $a2new
Replace-Text $a3old $a3new

Split-Pieces $a1new @($a1new, ".", $a2new, ".", $a3new)

Write-Host "STAGE-BODY-A OK"

# ---------------------------------------------------------------------------
# Main body paragraph, segment B (sentences 4-6). Sentence 6 splits into 3 runs.
# ---------------------------------------------------------------------------
$b1old = "As scientists venture into the intricate labyrinth of genetic information, they encounter a vast network of interactions between genes, proteins, and cellular processes"
$b1new = "From the smallest molecular components to the complex interactions of organ systems, Biology unveils a mesmerizing world of dynamic processes and remarkable adaptations"
$b2old = " Mutations, anomalies within the genetic sequence, can disrupt these interactions, leading to a kaleidoscope of medical conditions"
$b2new = " Delving into the realm of cells, we discover the fundamental building blocks of life, each possessing unique characteristics and functions"
$b3old = " By cracking the genetic code of disease, researchers embark on a quest to decipher the language of illness, revealing the root causes behind enigmatic symptoms"
$b3part1 = " Tissues, composed of specialized cells, form the framework of organs, which collectively execute intricate tasks that contribute to our overall physiology"
$b3part2 = "."
$b3part3 = " Marveling at the intricate interplay of these systems underscores the profound interconnectedness of the human body"
$b3new = $b3part1 + $b3part2 + $b3part3

Replace-Text $b1old $b1new
Replace-Text $b2old $b2new
Replace-Text $b3old $b3new

Split-Pieces $b1new @($b1new, ".", $b2new, ".", $b3part1, $b3part2, $b3part3)

Write-Host "STAGE-BODY-B OK"

# ---------------------------------------------------------------------------
# Main body paragraph, segment C (sentences 7-9)
# ---------------------------------------------------------------------------
$c1old = "Genome-wide association studies, the large-scale mapping of genetic variations, have revealed tantalizing clues linking DNA variations to disease susceptibility"
$c1new = "The study of Biology extends beyond the human body, inviting us to explore the vast panorama of life on Earth"
$c2old = " This emergent understanding underscores the influence of genetic factors in shaping our health destiny"
$c2new = " From microscopic organisms to towering trees, Biology delves into the diversity of species, examining their adaptations, interactions, and the intricate balance of ecosystems"
$c3old = " Yet, the odyssey continues, as researchers unravel the intricate tapestry of gene expression, the dynamic interplay between genes and the environment, and the dance of cellular pathways"
$c3new = " Engaging with Biology empowers us not only to understand ourselves better but also to appreciate the astonishing tapestry of life that surrounds us, inspiring awe and reverence for the natural world"

Replace-Text $c1old $c1new
Replace-Text $c2old $c2new
Replace-Text $c3old $c3new

Split-Pieces $c1new @($c1new, ".", $c2new, ".", $c3new)

Write-Host "STAGE-BODY-C OK"

Write-Host "DONE"
